# Applies the coinranking price/volume refresh described in the commit
# message "Updated cryptos list on Sat Apr 22 08:01:07 UTC 2023 with GitHub
# Actions" -- new Price (column D) / Volume(1h) (column E) snapshot values
# for the coin rows on Sheet1.
#
# D2:E51 are plain text cells (no numeric NumberFormat). Some of the new Price
# strings (e.g. "323.27") look like ordinary decimal numbers, so assigning them
# through .Value as-is would let Excel auto-convert the cell to a Number
# (re-typing/rounding the stored value). Each such value is therefore written
# with a leading apostrophe, exactly like a user typing `'323.27` into the cell
# would, to force Excel to keep treating it as text; the cell style is then put
# back to "Normal" so the forced quote-prefix does not leave a stray style
# behind (Excel tracks quote-prefixed cells via a cell style flag). Values that
# are not parseable as a plain number (e.g. "27.348.94", which has two dots)
# are left as-is since Excel stores those as text natively anyway.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "27.348.94"
    "E2" = "  -3.30%  "
    "D3" = "1.855.88"
    "E3" = "  -4.17%  "
    "E4" = "  -0.15%  "
    "D5" = "'323.27"
    "E5" = "  -2.58%  "
    "E6" = "  -0.10%  "
    "D7" = "'0.4521"
    "E7" = "  -4.58%  "
    "D8" = "'0.3865"
    "E8" = "  -5.08%  "
    "D9" = "'48.21"
    "E9" = "  -9.51%  "
    "D10" = "'0.07901"
    "E10" = "  -6.97%  "
    "E11" = "  -3.66%  "
    "D12" = "'21.33"
    "E12" = "  -4.67%  "
    "D13" = "1.863.20"
    "E13" = "  -2.42%  "
    "D14" = "'5.905"
    "E14" = "  -3.91%  "
    "D15" = "'7.108"
    "E15" = "  -6.12%  "
    "E16" = "  -0.24%  "
    "D17" = "'85.79"
    "E17" = "  -5.18%  "
    "D18" = "'0.00001029"
    "E18" = "  -3.93%  "
    "D19" = "'0.06548"
    "E19" = "  -0.52%  "
    "D20" = "'17.00"
    "E20" = "  -7.22%  "
    "E21" = "  -0.16%  "
    "D22" = "'5.539"
    "E22" = "  -4.49%  "
    "D23" = "27.331.07"
    "E23" = "  -3.34%  "
    "D24" = "'10.96"
    "E24" = "  -4.58%  "
    "D25" = "'2.281"
    "E25" = "  -0.47%  "
    "D26" = "2.081.04"
    "E26" = "  -2.80%  "
    "D27" = "'153.62"
    "D28" = "'19.90"
    "E28" = "  -1.42%  "
    "E29" = "  -5.04%  "
    "D30" = "'5.423"
    "E30" = "  -6.46%  "
    "D31" = "'120.90"
    "E31" = "  -2.40%  "
    "D32" = "'1.479"
    "E32" = "  +1.18%  "
    "D33" = "'0.09269"
    "E33" = "  -3.66%  "
    "D34" = "'0.9355"
    "E34" = "  -5.20%  "
    "D35" = "'3.600"
    "E35" = "  -1.38%  "
    "E36" = "  -6.09%  "
    "D37" = "'1.233"
    "E37" = "  -1.21%  "
    "E38" = "  -4.36%  "
    "D39" = "'0.05987"
    "E39" = "  -3.28%  "
    "D40" = "'8.168"
    "E40" = "  -11.90%  "
    "E41" = "  -0.16%  "
    "D42" = "'0.5897"
    "E42" = "  -5.08%  "
    "D43" = "'0.1896"
    "E43" = "  -0.89%  "
    "D44" = "'10.12"
    "E44" = "  -9.37%  "
    "D45" = "'1.274"
    "E45" = "  -3.03%  "
    "D46" = "'0.5618"
    "E46" = "  -5.19%  "
    "E47" = "  -6.30%  "
    "D48" = "'3.373"
    "E48" = "  -2.79%  "
    "D49" = "'1.917"
    "E49" = "  -6.75%  "
    "D50" = "'0.06765"
    "E50" = "  -0.45%  "
    "D51" = "'108.32"
    "E51" = "  -1.62%  "
}

foreach ($cellRef in $updates.Keys) {
    $value = $updates[$cellRef]
    $range = $ws.Range($cellRef)
    $range.Value = $value
    if ($value.StartsWith("'")) {
        # Drop the quote-prefix style Excel just added so only the cell
        # VALUE changes, matching a plain text cell with no explicit style.
        $range.Style = "Normal"
    }
}
